# Update Daily Report: 2026-02-05
# Appends a new day's block (date serial 46057) to Daily_Data, mirroring the
# previous day's block of 22 depository rows, then refreshes the rollups on
# Today_Summary and Monthly_Stats that are driven by that new data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Daily_Data: append the new day's block (rows 486:507) by copying the
#    prior day's block (rows 464:485), which carries over formatting/styles,
#    then correct the handful of cells that actually change day-over-day.
# ---------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily_Data")

$srcBlock = $daily.Range("A464:H485")
$dstBlock = $daily.Range("A486:H507")
$srcBlock.Copy($dstBlock)

# New block's Date column (A) advances one day: 46056 -> 46057
$daily.Range("A486:A507").Value = 46057

# BRINK'S, INC. Eligible (row 489): PREV_TOTAL (C) carries forward the prior
# day's TOTAL_TODAY (98856.745) rather than the prior day's own PREV_TOTAL
# that the straight copy brought over, and today's RECEIVED/NET_CHANGE (D/F)
# reset to 0 since the copied block brought over the prior day's activity.
$daily.Range("C489").Value = 98856.745
$daily.Range("D489").Value = 0
$daily.Range("F489").Value = 0

# MANFRA, TORDELLA & BROOKES, LLC Eligible (row 505): received 533.31 today.
$daily.Range("D505").Value = 533.3099999999999
$daily.Range("F505").Value = 533.3099999999999
$daily.Range("H505").Value = 1804.683

# ---------------------------------------------------------------------
# 2. Today_Summary: MANFRA, TORDELLA & BROOKES, LLC row (row 11) reflects
#    the new Eligible total and the resulting Total_Stock.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Today_Summary")
$summary.Range("B11").Value = 1804.683
$summary.Range("D11").Value = 52025.103

# ---------------------------------------------------------------------
# 3. Monthly_Stats: 2026-02 month-to-date rollup (row 2) and the
#    MANFRA, TORDELLA & BROOKES, LLC Eligible detail row (row 26) for 2026-02.
# ---------------------------------------------------------------------
$monthly = $wb.Worksheets.Item("Monthly_Stats")
$monthly.Range("B2").Value = 335980.783
$monthly.Range("D2").Value = 662618.7320000001

$monthly.Range("C26").Value = 533.3099999999999
$monthly.Range("E26").Value = 1804.683
